# Apply "ifoCAST full series evaluation" update to the qoq error table.
# Fills in the previously-missing staircase tail of each row (B2:K24)
# with the recomputed values for the full forecast horizon (Q0..Q9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 23,10
$data[0,0] = -0.9713005381550337
$data[0,1] = 0.2910874087102975
$data[0,2] = -0.1449385904233094
$data[0,3] = 0.5154458957241742
$data[0,4] = -0.1879827916858758
$data[0,5] = 0.04615233206547309
$data[0,6] = 0.08040507924932089
$data[0,7] = 0.9057460109392793
$data[0,8] = 0.249410584482785
$data[0,9] = -0.4089899291462971
$data[1,0] = 0.3128169072748576
$data[1,1] = 0.7245973247384453
$data[1,2] = -0.08695745817427486
$data[1,3] = 0.1004139362101289
$data[1,4] = 0.1144479125438225
$data[1,5] = 0.9310525395423103
$data[1,6] = 0.2709447666661666
$data[1,7] = -0.3890834678507667
$data[1,8] = 0.0092130989630414
$data[1,9] = -0.05258007222182071
$data[2,0] = -0.1030887313920102
$data[2,1] = 0.07647302587122451
$data[2,2] = 0.09012235911373839
$data[2,3] = 0.9090979817469648
$data[2,4] = 0.2510647147722038
$data[2,5] = -0.4077104680353097
$data[2,6] = -0.008762285969805494
$data[2,7] = -0.07024161732427531
$data[2,8] = -0.5044181462305073
$data[2,9] = 0.4743807131573582
$data[3,0] = 0.1720753300388297
$data[3,1] = 0.9182953219789133
$data[3,2] = 0.2313571522712326
$data[3,3] = -0.4369487835323326
$data[3,4] = -0.04090667078546828
$data[3,5] = -0.1032375837743004
$data[3,6] = -0.5376581911458389
$data[3,7] = 0.4410716177777917
$data[3,8] = 0.2139357511207785
$data[3,9] = -0.4195295179412606
$data[4,0] = 0.2025558615083408
$data[4,1] = -0.4390725672109059
$data[4,2] = -0.03362370600304576
$data[4,3] = -0.09282107083105617
$data[4,4] = -0.5262170457115171
$data[4,5] = 0.4528464178518407
$data[4,6] = 0.2258193719441708
$data[4,7] = -0.4076102730996941
$data[4,8] = 0.01979686535210479
$data[4,9] = -0.1465366454910707
$data[5,0] = -0.155505867098859
$data[5,1] = -0.1846662397643801
$data[5,2] = -0.6060897922829254
$data[5,3] = 0.3779170903031727
$data[5,4] = 0.1529367346250136
$data[5,5] = -0.4796464385677615
$data[5,6] = -0.0518896668375835
$data[5,7] = -0.2180789225979328
$data[5,8] = -0.09530727747379439
$data[5,9] = -0.04804865355017168
$data[6,0] = -0.5107750206255626
$data[6,1] = 0.4561847949904575
$data[6,2] = 0.2235272252108875
$data[6,3] = -0.412585116272613
$data[6,4] = 0.01354370435888463
$data[6,5] = -0.1533995520910376
$data[6,6] = -0.03097840982344718
$data[6,7] = 0.01611673723937751
$data[6,8] = -0.3805682542294263
$data[6,9] = -0.1265478981343244
$data[7,0] = 0.4836497670136274
$data[7,1] = -0.2716791846901883
$data[7,2] = 0.0989882439940924
$data[7,3] = -0.09374781479368632
$data[7,4] = 0.01667832048765883
$data[7,5] = 0.058195260878178
$data[7,6] = -0.3410837941741738
$data[7,7] = -0.0882697463915933
$data[7,8] = 0.1312785433800194
$data[7,9] = -0.1047779946421779
$data[8,0] = -0.09862000608775029
$data[8,1] = -0.2177076384565043
$data[8,2] = -0.07470130014448528
$data[8,3] = -0.018762957560511
$data[8,4] = -0.4116486652957548
$data[8,5] = -0.1559957098169479
$data[8,6] = 0.06481524009958373
$data[8,7] = -0.1706787749263282
$data[8,8] = -0.3495252635842865
$data[8,9] = -0.3600316430428294
$data[9,0] = -0.007299246851658558
$data[9,1] = 0.05611473536190403
$data[9,2] = -0.3334825685761003
$data[9,3] = -0.07638392106091824
$data[9,4] = 0.1450623997799764
$data[9,5] = -0.09015246042256481
$data[9,6] = -0.2688763367924936
$data[9,7] = -0.2793288770644347
$data[9,8] = -0.4870011301597575
$data[9,9] = -0.2697718993153602
$data[10,0] = -0.3458392694900738
$data[10,1] = -0.09057491098672893
$data[10,2] = 0.1301161455573347
$data[10,3] = -0.1054089187162663
$data[10,4] = -0.2842592456064009
$data[10,5] = -0.2947628657683481
$data[10,6] = -0.5024555220843632
$data[10,7] = -0.285234325920563
$data[10,8] = 0.02445755927687465
$data[10,9] = -0.1298439175827339
$data[11,0] = 0.3024188481467391
$data[11,1] = -0.01091678595725426
$data[11,2] = -0.2259338067192559
$data[11,3] = -0.2532411366013347
$data[11,4] = -0.4687411474418992
$data[11,5] = -0.2551474109997061
$data[11,6] = 0.05285908076345025
$data[11,7] = -0.1022254650337468
$data[11,8] = 0.5478383610817548
$data[11,9] = 0.3315048434816775
$data[12,0] = -0.3204423536017646
$data[12,1] = -0.2897555356160152
$data[12,2] = -0.4785434472588804
$data[12,3] = -0.2526553077156036
$data[12,4] = 0.06101006323956631
$data[12,5] = -0.09146971300380768
$data[12,6] = 0.5597931274032256
$data[12,7] = 0.3440115546243878
$data[12,8] = -0.06284176098887906
$data[12,9] = 0.4362810619427306
$data[13,0] = -0.3361387249989222
$data[13,1] = -0.172079154149077
$data[13,2] = 0.1140244002424372
$data[13,3] = -0.05074182916399023
$data[13,4] = 0.5950432373871204
$data[13,5] = 0.376819120134967
$data[13,6] = -0.03112348437038687
$data[13,7] = 0.4675134828696286
$data[13,8] = 0.308484786548516
$data[13,9] = 0.3645688493654578
$data[14,0] = 0.2919567931684033
$data[14,1] = 0.03611576641084058
$data[14,2] = 0.6379689631723193
$data[14,3] = 0.3987277948059625
$data[14,4] = -0.01926926647574417
$data[14,5] = 0.4745601165931918
$data[14,6] = 0.3132336451036749
$data[14,7] = 0.3682199310910328
$data[14,8] = 2.617463111368334
$data[14,9] = 10.0863474470963
$data[15,0] = 0.04668378654349858
$data[15,1] = 0.646743751744091
$data[15,2] = 0.4065301988496979
$data[15,3] = -0.01195920564123576
$data[15,4] = 0.4816209642468198
$data[15,5] = 0.3201688216047165
$data[15,6] = 0.3750919343242142
$data[15,7] = 2.624303448075727
$data[15,8] = 10.09317195097463
$data[15,9] = -8.078224169346534
$data[16,0] = 0.5397693479284105
$data[16,1] = 0.3412705759114291
$data[16,2] = -0.05754980417421507
$data[16,3] = 0.4452902877305661
$data[16,4] = 0.2881987268708067
$data[16,5] = 0.3451757067707001
$data[16,6] = 2.595354809876611
$data[16,7] = 10.06467924472723
$data[16,8] = -8.106501994311596
$data[16,9] = 0.03320682975976369
$data[17,0] = 0.3751437736559251
$data[17,1] = -0.0547966439174179
$data[17,2] = 0.4337013568798115
$data[17,3] = 0.2700773128138926
$data[17,4] = 0.3240755934079451
$data[17,5] = 2.572896473190138
$data[17,6] = 10.04160153156291
$data[17,7] = -8.129862178460501
$data[17,8] = 0.009717811984482905
$data[17,9] = 2.157246767248676
$data[18,0] = -0.312072073423668
$data[18,1] = 0.3024220853976149
$data[18,2] = 0.1967216315401449
$data[18,3] = 0.2768343187407715
$data[18,4] = 2.5374495888114
$data[18,5] = 10.01148064496908
$data[18,6] = -8.15757796247812
$data[18,7] = -0.01691188077547812
$data[18,8] = 2.13110753017287
$data[18,9] = -1.308719802433826
$data[19,0] = 0.259157971428892
$data[19,1] = 0.16838023551327
$data[19,2] = 0.2521582932286491
$data[19,3] = 2.513162955834568
$data[19,4] = 9.986857348713633
$data[19,5] = -8.182561135418815
$data[19,6] = -0.04214184030638374
$data[19,7] = 2.105731958361166
$data[19,8] = -1.334174982466894
$data[19,9] = -1.395299246927819
$data[20,0] = 0.2785713907394387
$data[20,1] = 0.3082809271729959
$data[20,2] = 2.547244771768479
$data[20,3] = 10.01195673258608
$data[20,4] = -8.161122484703073
$data[20,5] = -0.02219513091605813
$data[20,6] = 2.12507060979506
$data[20,7] = -1.315084157547472
$data[20,8] = -1.376309431020366
$data[20,9] = 0.665125583809529
$data[21,0] = 0.1514308419078968
$data[21,1] = 2.452077041002839
$data[21,2] = 9.957529601920154
$data[21,3] = -8.197676472114839
$data[21,4] = -0.04986031939925989
$data[21,5] = 2.101663470123121
$data[21,6] = -1.336421077423499
$data[21,7] = -1.396642028765399
$data[21,8] = 0.6452816047775823
$data[21,9] = 0.06291644308519029
$data[22,0] = 2.403213427390482
$data[22,1] = 9.929760751643483
$data[22,2] = -8.214896317152331
$data[22,3] = -0.06195543535285469
$data[22,4] = 2.092079342639022
$data[22,5] = -1.344778008240451
$data[22,6] = -1.404398746163505
$data[22,7] = 0.6378183817016827
$data[22,8] = 0.05559674216991822
$data[22,9] = 0.1514544743675568

$ws.Range("B2:K24").Value2 = $data
